# "code to write" is the active / tab-selected sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill colors used by this tracker sheet (OLE BGR-packed RGB longs):
#   green  FF92D050 -> "Yes"          (already used for C2:C5 before the edit)
#   yellow FFFFFF00 -> "In Progress"  (already used for C7:C10 before the edit)
$green  = 5296274   # RGB(0x92, 0xD0, 0x50)
$yellow = 65535      # RGB(0xFF, 0xFF, 0x00)

# C2 ("MdlMain" / Main) moves from done ("Yes") to "In Progress".
$ws.Range("C2").Value = "In Progress"
$ws.Range("C2").Interior.Color = $yellow

# Rows that were "N" (no status fill) are now marked done ("Yes" / green).
$doneRows = @(11, 12, 13, 14, 16, 18)
foreach ($r in $doneRows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = "Yes"
    $cell.Interior.Color = $green
}

# Row 17 ("Server" / CleanActivateJobWebParams) previously had no status cell
# at all; it now gets "Yes" (green) as well.
$ws.Range("C17").Value = "Yes"
$ws.Range("C17").Interior.Color = $green

# Row 15 ("MdlServer" / fCheckForDuplicateRealTimes) previously had no status
# cell; it now gets "In Progress" (yellow).
$ws.Range("C15").Value = "In Progress"
$ws.Range("C15").Interior.Color = $yellow

# Row 23 ("Machine" / Init) moves from "N" (no fill) to "In Progress" (yellow).
$ws.Range("C23").Value = "In Progress"
$ws.Range("C23").Interior.Color = $yellow

# Update the sheet's last active cell/selection to reflect where work left off.
$ws.Range("C23").Select() | Out-Null
